# Update the "Förändrad" (Changed) date column C for every data row (2-396)
# from 45184 (2023-09-15) to 45186 (2023-09-17).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C396").Value2 = 45186

# For the first 14 data rows (2-15), which contain populated species/hyperlink
# columns, append a second (friendly-name) argument to each HYPERLINK formula
# in columns S, T, V, W, X, Y. The friendly name is the designation value
# found in column A of the same row (e.g. "A 59085-2018").
$hyperlinkCols = @("S", "T", "V", "W", "X", "Y")

for ($row = 2; $row -le 15; $row++) {
    $designation = $ws.Range("A" + $row).Value2

    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Range($col + $row)
        $formula = $cell.Formula

        if ($formula -and $formula.EndsWith(')')) {
            $newFormula = $formula.Substring(0, $formula.Length - 1) + ', "' + $designation + '")'
            $cell.Formula = $newFormula
        }
    }
}
